$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (shared strings with rich-text runs) ---

# "Volume 32   Number  31" -> "...Number  32" (report issue number)
$volCell = $ws.Range("A8")
$volText = $volCell.Value()
$volIdx = $volText.LastIndexOf("31")
$volChars = $volCell.Characters($volIdx + 1, 2)
$volChars.Text = "32"

# "Report Covering the Week  7/28/2025  Through  8/3/2025"
#   -> "...Week  8/4/2025  Through  8/10/2025"
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value()
$startIdx = $weekText.IndexOf("7/28/2025")
$startChars = $weekCell.Characters($startIdx + 1, 9)
$startChars.Text = "8/4/2025"

$weekText2 = $weekCell.Value()
$endIdx = $weekText2.IndexOf("8/3/2025")
$endChars = $weekCell.Characters($endIdx + 1, 8)
$endChars.Text = "8/10/2025"

# --- Crime-statistics table updates (rows 14-30) ---

$ws.Range("J14").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("J14").Copy($ws.Range("F14"))
$ws.Range("F14").Value = 1
$ws.Range("J14").Copy($ws.Range("I14"))
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -75
$ws.Range("J14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("J14").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -28.571428571428
$ws.Range("D16").Value = 1
$ws.Range("G16").Value = 7
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = -32.5
$ws.Range("L16").Value = 17.391304347826
$ws.Range("M16").Value = -28.947368421052
$ws.Range("N16").Value = -83.536585365853
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 70
$ws.Range("J17").Value = 105
$ws.Range("K17").Value = -33.333333333333
$ws.Range("L17").Value = -1.408450704225
$ws.Range("M17").Value = 29.629629629629
$ws.Range("N17").Value = -45.736434108527
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = -5.882352941176
$ws.Range("L18").Value = 39.130434782608
$ws.Range("M18").Value = -11.111111111111
$ws.Range("N18").Value = -88.363636363636
$ws.Range("D14").Copy($ws.Range("C19"))
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = -100
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = -42.857142857142
$ws.Range("J19").Value = 95
$ws.Range("K19").Value = -22.105263157894
$ws.Range("L19").Value = -36.752136752136
$ws.Range("M19").Value = 4.225352112676
$ws.Range("N19").Value = -39.344262295082
$ws.Range("D14").Copy($ws.Range("C20"))
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -40
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -54.285714285714
$ws.Range("L20").Value = -33.333333333333
$ws.Range("N20").Value = -92.825112107623
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 34
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = -22.727272727272
$ws.Range("I21").Value = 225
$ws.Range("J21").Value = 312
$ws.Range("K21").Value = -27.884615384615
$ws.Range("L21").Value = -13.461538461538
$ws.Range("M21").Value = 2.739726027397
$ws.Range("N21").Value = -75.649350649350
$ws.Range("D14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("J14").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 2
$ws.Range("D14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 13
$ws.Range("K23").Value = -31.578947368421
$ws.Range("L23").Value = 62.5
$ws.Range("M23").Value = -23.529411764705
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 50
$ws.Range("H24").Value = -18
$ws.Range("I24").Value = 309
$ws.Range("J24").Value = 281
$ws.Range("K24").Value = 9.964412811387
$ws.Range("L24").Value = 6.920415224913
$ws.Range("M24").Value = 75.568181818181
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -19.230769230769
$ws.Range("I25").Value = 156
$ws.Range("J25").Value = 138
$ws.Range("K25").Value = 13.043478260869
$ws.Range("L25").Value = 4
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 400
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -8.333333333333
$ws.Range("I26").Value = 144
$ws.Range("J26").Value = 166
$ws.Range("K26").Value = -13.253012048192
$ws.Range("L26").Value = -1.369863013698
$ws.Range("M26").Value = -33.944954128440
$ws.Range("J14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 50
$ws.Range("J14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 13
$ws.Range("K28").Value = 30
$ws.Range("L28").Value = -23.529411764705
$ws.Range("J14").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 2
$ws.Range("J14").Copy($ws.Range("F29"))
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = -42.857142857142
$ws.Range("L29").Value = -55.555555555555
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -20
$ws.Range("J14").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 2
$ws.Range("J14").Copy($ws.Range("F30"))
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = -20
$ws.Range("L30").Value = -20
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
